# "Generate Report for Handback" — the d8525287-ace8-4478-9dad-43996428ca5f
# item has come back from handback for both zh-cn and de-de, so:
#   * the Overview sheet's Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   * each language sheet's Status does the same, and the previously-empty
#     "Latest Target File" / "Latest Handback File" columns get filled in
#     (as real hyperlinks, matching the pattern already used by the other
#     rows), and "Latest Handback DateTime" moves off the zero-date
#     placeholder to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B5").Value = $statusHandedBack
$ov.Range("C5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C5").Value = $statusHandedBack

$zh.Hyperlinks.Add(
    $zh.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/d8525287ace84478/e2e/d8525287-ace8-4478-9dad-43996428ca5f.md",
    $null,
    $null,
    "d8525287-ace8-4478-9dad-43996428ca5f.md"
)
$zh.Range("F5").Font.Underline = $true
$zh.Range("F5").Font.Color = 15570276
$zh.Range("F5").Font.Name = "Calibri"
$zh.Range("F5").Font.Size = 11

$zh.Hyperlinks.Add(
    $zh.Range("G5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d8525287ace84478/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/d8525287-ace8-4478-9dad-43996428ca5f.4f3a3ce3b866c63a7aa3c428b2f49b69ee48eb6f.zh-cn.xlf",
    $null,
    $null,
    "d8525287-ace8-4478-9dad-43996428ca5f.4f3a3ce3b866c63a7aa3c428b2f49b69ee48eb6f.zh-cn.xlf"
)
$zh.Range("G5").Font.Underline = $true
$zh.Range("G5").Font.Color = 15570276
$zh.Range("G5").Font.Name = "Calibri"
$zh.Range("G5").Font.Size = 11

$zh.Range("H5").Value = "2016-03-30 10:19:48"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C5").Value = $statusHandedBack

$de.Hyperlinks.Add(
    $de.Range("F5"),
    "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/d8525287ace84478/e2e/d8525287-ace8-4478-9dad-43996428ca5f.md",
    $null,
    $null,
    "d8525287-ace8-4478-9dad-43996428ca5f.md"
)
$de.Range("F5").Font.Underline = $true
$de.Range("F5").Font.Color = 15570276
$de.Range("F5").Font.Name = "Calibri"
$de.Range("F5").Font.Size = 11

$de.Hyperlinks.Add(
    $de.Range("G5"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d8525287ace84478/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/d8525287-ace8-4478-9dad-43996428ca5f.4f3a3ce3b866c63a7aa3c428b2f49b69ee48eb6f.de-de.xlf",
    $null,
    $null,
    "d8525287-ace8-4478-9dad-43996428ca5f.4f3a3ce3b866c63a7aa3c428b2f49b69ee48eb6f.de-de.xlf"
)
$de.Range("G5").Font.Underline = $true
$de.Range("G5").Font.Color = 15570276
$de.Range("G5").Font.Name = "Calibri"
$de.Range("G5").Font.Size = 11

$de.Range("H5").Value = "2016-03-30 10:20:06"
